$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("James Harden",          "PG,SG",    "LA Clippers"),
  @("Anfernee Simons",       "PG,SG",    "Portland Trail Blazers"),
  @("Collin Sexton",         "PG,SG",    "Utah Jazz"),
  @("Rui Hachimura",         "SF,PF",    "Los Angeles Lakers"),
  @("Paul George",           "SG,SF,PF", "Philadelphia 76ers"),
  @("Anthony Edwards",       "SG,SF",    "Minnesota Timberwolves"),
  @("Aaron Gordon",          "PF,C",     "Denver Nuggets"),
  @("Jaren Jackson Jr.",     "PF,C",     "Memphis Grizzlies"),
  @("Ivica Zubac",           "C",        "LA Clippers"),
  @("Tyus Jones",            "PG",       "Phoenix Suns"),
  @("Fred VanVleet",         "PG",       "Houston Rockets"),
  @("Khris Middleton",       "SF",       "Milwaukee Bucks"),
  @("Jaden Ivey",            "PG,SG",    "Detroit Pistons"),
  @("Giannis Antetokounmpo", "PF,C",     "Milwaukee Bucks"),
  @("Jayson Tatum",          "SF,PF",    "Boston Celtics"),
  @("Zion Williamson",       "PF,C",     "New Orleans Pelicans"),
  @("Kyle Kuzma",            "PF",       "Washington Wizards"),
  @("Andrew Wiggins",        "SF,PF",    "Golden State Warriors")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
}
